$d = $word.ActiveDocument

# --- 1. Split the last paragraph's "normal" text into two runs: "N" and "ormal" ---
# (same, empty run formatting on both -- achieved by doing the edit as a tracked
#  change and then accepting it, which preserves the run boundary that a plain
#  in-place text edit would otherwise coalesce away)
$lastTextPara = $d.Paragraphs($d.Paragraphs.Count)
$fullRange = $lastTextPara.Range
$paraStart = $fullRange.Start

$d.TrackRevisions = $true
$firstChar = $d.Range($paraStart, $paraStart + 1)
$firstChar.Text = ""
$insertPoint = $d.Range($paraStart, $paraStart)
$insertPoint.InsertBefore("N")
$d.TrackRevisions = $false
$d.AcceptAllRevisions() | Out-Null

# --- 2. Add a new paragraph after it with mixed (bold+italic+underline) formatting
#         in the middle of a sentence ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newRange = $newPara.Range
$newStart = $newRange.Start

$sentence = "This is a sentence with all formatting options in the middle of the sentence."
$newRange.Text = $sentence

$allStart = $newStart + $sentence.IndexOf("all")
$allEnd = $allStart + 3
$allRange = $d.Range($allStart, $allEnd)
$allRange.Font.Bold = 1
$allRange.Font.Italic = 1
$allRange.Font.Underline = 1

# --- 3. styles.xml: mark the built-in "Normal Table" style as a quick style ---
$tableNormal = $d.Styles("Normal Table")
$tableNormal.QuickStyle = $true

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
